# Updated cryptos list on Fri Mar 15 16:29:18 UTC 2024 with GitHub Actions
# Refresh Price / Volume(1h) figures for the coinranking.com snapshot, and
# re-sort a couple of rows whose rank order changed (Uniswap/Polygon and
# ApeXProtocol/Stellar/WEMIXToken).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.796.52'
$ws.Range("E2").Value = '  -3.82%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.675.40'
$ws.Range("E3").Value = '  -4.55%  '

$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.83'
$ws.Range("E5").Value = '  +0.31%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '179.40'
$ws.Range("E6").Value = '  +8.16%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.668.15'
$ws.Range("E7").Value = '  -4.50%  '

$ws.Range("E8").Value = '  -6.58%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.01'
$ws.Range("E9").Value = '  +0.60%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.712'
$ws.Range("E10").Value = '  -4.24%  '

$ws.Range("E11").Value = '  -6.97%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '55.31'
$ws.Range("E12").Value = '  +4.45%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000287'
$ws.Range("E13").Value = '  -10.09%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.34'
$ws.Range("E14").Value = '  -7.94%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.250.87'
$ws.Range("E15").Value = '  -4.93%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.675.86'
$ws.Range("E16").Value = '  -5.11%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.29'
$ws.Range("E17").Value = '  -6.29%  '

$ws.Range("E18").Value = '  -2.28%  '

# Row 19/20 swap: Uniswap now ranks above Polygon
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.72'
$ws.Range("E19").Value = '  -7.34%  '

$ws.Range("B20").Value = 'Polygon'
$ws.Range("C20").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.12'
$ws.Range("E20").Value = '  -6.80%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '67.577.65'
$ws.Range("E21").Value = '  -4.02%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '407.33'
$ws.Range("E22").Value = '  -6.08%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.52'
$ws.Range("E23").Value = '  -3.92%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '87.94'
$ws.Range("E24").Value = '  -6.36%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.99'
$ws.Range("E25").Value = '  -8.20%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.72'
$ws.Range("E26").Value = '  -7.16%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.86'
$ws.Range("E27").Value = '  -0.55%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.84'
$ws.Range("E28").Value = '  -5.40%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.06'
$ws.Range("E29").Value = '  +2.20%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.39'
$ws.Range("E30").Value = '  -8.17%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.46'
$ws.Range("E31").Value = '  -7.01%  '

$ws.Range("E32").Value = '  -8.07%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '12.41'
$ws.Range("E33").Value = '  -7.45%  '

$ws.Range("E34").Value = '  -6.69%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '42.92'
$ws.Range("E35").Value = '  -11.21%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '63.94'
$ws.Range("E36").Value = '  -7.62%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '591.15'
$ws.Range("E37").Value = '  -3.80%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0869'
$ws.Range("E38").Value = '  -9.93%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  +0.09%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.397'
$ws.Range("E40").Value = '  -4.90%  '

$ws.Range("E41").Value = '  -0.11%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.135'
$ws.Range("E42").Value = '  -4.73%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.71'
$ws.Range("E43").Value = '  +1.50%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.98'
$ws.Range("E44").Value = '  -8.04%  '

$ws.Range("E45").Value = '  -7.11%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.86'
$ws.Range("E46").Value = '  -10.48%  '

$ws.Range("E47").Value = '  -8.09%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.760.34'
$ws.Range("E48").Value = '  -2.42%  '

# Rows 49-51 re-sort: ApeXProtocol, Stellar, WEMIXToken (cyclic shift)
$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.16'
$ws.Range("E49").Value = '  -3.86%  '

$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.133'
$ws.Range("E50").Value = '  -6.74%  '

$ws.Range("B51").Value = 'WEMIXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.68'
$ws.Range("E51").Value = '  -3.55%  '
